$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated GWL (global warming level) period values for columns B:E (rows 2-35),
# recalculated with the new CMIP6 ensemble. Column A (ensemble member names)
# is unchanged.
$values = @(
    @("B2", ""),
    @("C2", ""),
    @("D2", ""),
    @("E2", ""),
    @("B3", "2040-2059"),
    @("C3", ""),
    @("D3", ""),
    @("E3", ""),
    @("B4", "2027-2046"),
    @("C4", ""),
    @("D4", ""),
    @("E4", ""),
    @("B5", ""),
    @("C5", ""),
    @("D5", ""),
    @("E5", ""),
    @("B6", "2011-2030"),
    @("C6", "2044-2063"),
    @("D6", ""),
    @("E6", ""),
    @("B7", "2029-2048"),
    @("C7", ""),
    @("D7", ""),
    @("E7", ""),
    @("B8", "2016-2035"),
    @("C8", "2031-2050"),
    @("D8", ""),
    @("E8", ""),
    @("B9", ""),
    @("C9", ""),
    @("D9", ""),
    @("E9", ""),
    @("B10", "2045-2064"),
    @("C10", "2057-2076"),
    @("D10", ""),
    @("E10", ""),
    @("B11", "2039-2058"),
    @("C11", "2047-2066"),
    @("D11", ""),
    @("E11", ""),
    @("B12", "2036-2055"),
    @("C12", "2047-2066"),
    @("D12", ""),
    @("E12", ""),
    @("B13", "2039-2058"),
    @("C13", "2050-2069"),
    @("D13", ""),
    @("E13", ""),
    @("B14", "2022-2041"),
    @("C14", "2042-2061"),
    @("D14", ""),
    @("E14", ""),
    @("B15", "2021-2040"),
    @("C15", "2055-2074"),
    @("D15", ""),
    @("E15", ""),
    @("B16", "2028-2047"),
    @("C16", "2062-2081"),
    @("D16", ""),
    @("E16", ""),
    @("B17", "2028-2047"),
    @("C17", "2060-2079"),
    @("D17", ""),
    @("E17", ""),
    @("B18", "2027-2046"),
    @("C18", "2037-2056"),
    @("D18", ""),
    @("E18", ""),
    @("B19", "2015-2034"),
    @("C19", "2030-2049"),
    @("D19", "2054-2073"),
    @("E19", ""),
    @("B20", "2017-2036"),
    @("C20", "2033-2052"),
    @("D20", "2060-2079"),
    @("E20", ""),
    @("B21", "2029-2048"),
    @("C21", "2078-2097"),
    @("D21", ""),
    @("E21", ""),
    @("B22", "2023-2042"),
    @("C22", "2060-2079"),
    @("D22", ""),
    @("E22", ""),
    @("B23", "2032-2051"),
    @("C23", "2045-2064"),
    @("D23", "2064-2083"),
    @("E23", ""),
    @("B24", "2027-2046"),
    @("C24", "2037-2056"),
    @("D24", "2056-2075"),
    @("E24", "2072-2091"),
    @("B25", "2022-2041"),
    @("C25", "2035-2054"),
    @("D25", "2051-2070"),
    @("E25", "2067-2086"),
    @("B26", "2031-2050"),
    @("C26", "2043-2062"),
    @("D26", "2066-2085"),
    @("E26", "2079-2098"),
    @("B27", "2018-2037"),
    @("C27", "2030-2049"),
    @("D27", "2049-2068"),
    @("E27", "2066-2085"),
    @("B28", "2020-2039"),
    @("C28", "2048-2067"),
    @("D28", "2063-2082"),
    @("E28", "2073-2092"),
    @("B29", "2019-2038"),
    @("C29", "2042-2061"),
    @("D29", "2067-2086"),
    @("E29", ""),
    @("B30", "2041-2060"),
    @("C30", "2048-2067"),
    @("D30", "2062-2081"),
    @("E30", ""),
    @("B31", "2025-2044"),
    @("C31", "2032-2051"),
    @("D31", "2047-2066"),
    @("E31", "2058-2077"),
    @("B32", "2009-2028"),
    @("C32", "2023-2042"),
    @("D32", "2045-2064"),
    @("E32", "2056-2075"),
    @("B33", "2009-2028"),
    @("C33", "2024-2043"),
    @("D33", "2046-2065"),
    @("E33", "2057-2076"),
    @("B34", "2034-2053"),
    @("C34", "2041-2060"),
    @("D34", "2066-2085"),
    @("E34", "2081-2100"),
    @("B35", "2028-2047"),
    @("C35", "2039-2058"),
    @("D35", "2057-2076"),
    @("E35", "2068-2087")
)

foreach ($pair in $values) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# Restore the last active cell selection recorded for the sheet.
$null = $ws.Range("I21").Select()

